$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gof")
$ws.Range("D2").Value = 26042
$ws.Range("F2").Value = 26096
$ws.Range("G2").Value = 26239
$ws.Range("D3").Value = 26022
$ws.Range("F3").Value = 26108
$ws.Range("G3").Value = 26336

$ws = $wb.Worksheets.Item("facets")
$ws.Range("B2").Value = 768
$ws.Range("B3").Value = 732

$ws = $wb.Worksheets.Item("Estimates 0-1")
$ws.Range("B2").Value = 0.031
$ws.Range("D2").Value = 0.025
$ws.Range("E2").Value = 0.111
$ws.Range("B3").Value = -0.053
$ws.Range("C3").Value = 0.089
$ws.Range("D3").Value = -0.042
$ws.Range("E3").Value = 0.355
$ws.Range("B4").Value = -0.064
$ws.Range("D4").Value = -0.051
$ws.Range("E4").Value = 1.066
$ws.Range("B5").Value = -0.055
$ws.Range("C5").Value = 0.083
$ws.Range("D5").Value = -0.044
$ws.Range("E5").Value = 0.439
$ws.Range("B6").Value = -0.017
$ws.Range("C6").Value = 0.081
$ws.Range("D6").Value = -0.014
$ws.Range("E6").Value = 0.044
$ws.Range("B7").Value = -0.03
$ws.Range("C7").Value = 0.08
$ws.Range("D7").Value = -0.024
$ws.Range("E7").Value = 0.141
$ws.Range("B8").Value = 0.052
$ws.Range("D8").Value = 0.042
$ws.Range("E8").Value = 0.927
$ws.Range("B9").Value = 0.247
$ws.Range("D9").Value = 0.198
$ws.Range("E9").Value = 10.028
$ws.Range("B10").Value = 0.075
$ws.Range("C10").Value = 0.078
$ws.Range("D10").Value = 0.06
$ws.Range("E10").Value = 0.925
$ws.Range("B11").Value = 0.062
$ws.Range("C11").Value = 0.079
$ws.Range("D11").Value = 0.05
$ws.Range("E11").Value = 0.616
$ws.Range("B12").Value = 0.186
$ws.Range("C12").Value = 0.081
$ws.Range("D12").Value = 0.149
$ws.Range("E12").Value = 5.273
$ws.Range("B13").Value = 0.047
$ws.Range("C13").Value = 0.084
$ws.Range("D13").Value = 0.038
$ws.Range("E13").Value = 0.313
$ws.Range("B14").Value = -0.103
$ws.Range("C14").Value = 0.086
$ws.Range("D14").Value = -0.082
$ws.Range("E14").Value = 1.434
$ws.Range("B15").Value = -0.224
$ws.Range("C15").Value = 0.051
$ws.Range("D15").Value = -0.179
$ws.Range("E15").Value = 19.291
$ws.Range("B16").Value = -0.014
$ws.Range("C16").Value = 0.101
$ws.Range("D16").Value = -0.011
$ws.Range("E16").Value = 0.019
$ws.Range("B17").Value = -0.179
$ws.Range("D17").Value = -0.143
$ws.Range("E17").Value = 15.142
$ws.Range("B18").Value = -0.041
$ws.Range("C18").Value = 0.312
$ws.Range("D18").Value = -0.033
$ws.Range("E18").Value = 0.017

$ws = $wb.Worksheets.Item("Main effect 0-1")
$ws.Range("B2").Value = -0.34
$ws.Range("C2").Value = -0.272
$ws.Range("B3").Value = -0.242
$ws.Range("C3").Value = -0.194
